$wb = $excel.ActiveWorkbook

# Update "想去人数" (interested count) values that changed upstream.
# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 314
$ws1.Range("F4").Value = 49

# Sheet "全部类型" (all types) mirrors the same data
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 314
$ws4.Range("F4").Value = 49
